$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 7183.625
$ws.Range("J51").Value = 7098.75
$ws.Range("L51").Value = 7098.75
$ws.Range("N51").Value = -8066.75
$ws.Range("H112").Value = 2521.3
$ws.Range("I112").Value = 1564.3334
$ws.Range("J112").Value = 2931.4285
$ws.Range("K112").Value = 4693.0002
$ws.Range("L112").Value = 8794.2855
$ws.Range("M112").Value = -3585.0002
$ws.Range("N112").Value = -11010.2855
$ws.Range("H129").Value = 2681
$ws.Range("I129").Value = 1190
$ws.Range("J129").Value = 3053.75
$ws.Range("K129").Value = 3570
$ws.Range("L129").Value = 9161.25
$ws.Range("M129").Value = 1430
$ws.Range("N129").Value = -19161.25
$ws.Range("H131").Value = 4017.625
$ws.Range("J131").Value = 7900
$ws.Range("L131").Value = 23700
$ws.Range("N131").Value = -33780
$ws.Range("H138").Value = 2221.7222
$ws.Range("J138").Value = 3483.5
$ws.Range("L138").Value = 10450.5
$ws.Range("N138").Value = -20730.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 95962
$ws.Range("I45").Value = 187170.9
$ws.Range("J45").Value = 4753.091
$ws.Range("K45").Value = 187170.9
$ws.Range("L45").Value = 4753.091
$ws.Range("M45").Value = -186793.9
$ws.Range("N45").Value = -5507.091
$ws.Range("H97").Value = 14293729
$ws.Range("I97").Value = 13164
$ws.Range("K97").Value = 13164
$ws.Range("M97").Value = -12668
$ws.Range("H132").Value = 3443.25
$ws.Range("I132").Value = 3508
$ws.Range("J132").Value = 3335.3333
$ws.Range("K132").Value = 10524
$ws.Range("L132").Value = 10005.9999
$ws.Range("M132").Value = -7994
$ws.Range("N132").Value = -15065.9999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3241.95
$ws.Range("J20").Value = 4354.1
$ws.Range("L20").Value = 4354.1
$ws.Range("N20").Value = -4848.1
$ws.Range("H86").Value = 5350.273
$ws.Range("I86").Value = 6935.7334
$ws.Range("J86").Value = 1952.8572
$ws.Range("K86").Value = 6935.7334
$ws.Range("L86").Value = 1952.8572
$ws.Range("M86").Value = -5812.7334
$ws.Range("N86").Value = -4198.8572
$ws.Range("H89").Value = 5350.273
$ws.Range("I89").Value = 6935.7334
$ws.Range("J89").Value = 1952.8572
$ws.Range("K89").Value = 34678.667
$ws.Range("L89").Value = 9764.286
$ws.Range("M89").Value = -29062.667
$ws.Range("N89").Value = -20996.286
$ws.Range("H99").Value = 15010.261
$ws.Range("I99").Value = 18339.438
$ws.Range("J99").Value = 7400.7144
$ws.Range("K99").Value = 18339.438
$ws.Range("L99").Value = 7400.7144
$ws.Range("M99").Value = -16841.438
$ws.Range("N99").Value = -10396.7144
$ws.Range("H105").Value = 96342.37
$ws.Range("I105").Value = 250975
$ws.Range("K105").Value = 250975
$ws.Range("M105").Value = -249228
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1673.8
$ws.Range("J16").Value = 2170.75
$ws.Range("L16").Value = 2170.75
$ws.Range("N16").Value = -2744.75
$ws.Range("H31").Value = 12575.571
$ws.Range("I31").Value = 18051.125
$ws.Range("J31").Value = 5274.8335
$ws.Range("K31").Value = 18051.125
$ws.Range("L31").Value = 5274.8335
$ws.Range("M31").Value = -17756.125
$ws.Range("N31").Value = -5864.8335
$ws.Range("H34").Value = 12575.571
$ws.Range("I34").Value = 18051.125
$ws.Range("J34").Value = 5274.8335
$ws.Range("K34").Value = 18051.125
$ws.Range("L34").Value = 5274.8335
$ws.Range("M34").Value = -17849.125
$ws.Range("N34").Value = -5678.8335
$ws.Range("H99").Value = 12504449
$ws.Range("I99").Value = 20836264
$ws.Range("K99").Value = 20836264
$ws.Range("M99").Value = -20834766
$ws.Range("H107").Value = 4275.5884
$ws.Range("I107").Value = 5803.696
$ws.Range("J107").Value = 1080.4546
$ws.Range("K107").Value = 5803.696
$ws.Range("L107").Value = 1080.4546
$ws.Range("M107").Value = -3883.696
$ws.Range("N107").Value = -4920.4546
$ws.Range("H113").Value = 1673.8
$ws.Range("J113").Value = 2170.75
$ws.Range("L113").Value = 2170.75
$ws.Range("N113").Value = -6510.75
$ws.Range("H122").Value = 11549.5
$ws.Range("I122").Value = 16335.125
$ws.Range("J122").Value = 1978.25
$ws.Range("K122").Value = 49005.375
$ws.Range("L122").Value = 5934.75
$ws.Range("M122").Value = -46555.375
$ws.Range("N122").Value = -10834.75
$ws.Range("H126").Value = 12504449
$ws.Range("I126").Value = 20836264
$ws.Range("K126").Value = 62508792
$ws.Range("M126").Value = -62506322
$ws.Range("H132").Value = 1923.2222
$ws.Range("I132").Value = 2118.5
$ws.Range("J132").Value = 1532.6666
$ws.Range("K132").Value = 6355.5
$ws.Range("L132").Value = 4597.9998
$ws.Range("M132").Value = -3825.5
$ws.Range("N132").Value = -9657.9998
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 386
$ws.Range("I92").Value = 279.5
$ws.Range("J92").Value = 492.5
$ws.Range("K92").Value = 838.5
$ws.Range("L92").Value = 1477.5
$ws.Range("M92").Value = 409.5
$ws.Range("N92").Value = -3973.5
$ws.Range("H95").Value = 4450
$ws.Range("I95").Value = 4400
$ws.Range("J95").Value = 4500
$ws.Range("K95").Value = 13200
$ws.Range("L95").Value = 13500
$ws.Range("M95").Value = -11141
$ws.Range("N95").Value = -17618
$ws.Range("H101").Value = 4999
$ws.Range("J101").Value = 4999
$ws.Range("L101").Value = 14997
$ws.Range("N101").Value = -19865
$ws.Range("H105").Value = 9999.878000000001
$ws.Range("J105").Value = 9999.878000000001
$ws.Range("L105").Value = 29999.634
$ws.Range("N105").Value = -35241.63400000001
$ws.Range("H106").Value = 6990
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("H131").Value = 1462.11
$ws.Range("I131").Value = 830.53845
$ws.Range("J131").Value = 1556.4828
$ws.Range("K131").Value = 2491.61535
$ws.Range("L131").Value = 4669.4484
$ws.Range("M131").Value = 2548.38465
$ws.Range("N131").Value = -14749.4484
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 344000
$ws.Range("J95").Value = 344000
$ws.Range("L95").Value = 344000
$ws.Range("N95").Value = -349492
$ws.Range("H113").Value = 12355.454
$ws.Range("J113").Value = 3922.5
$ws.Range("L113").Value = 3922.5
$ws.Range("N113").Value = -8262.5
$ws.Range("H126").Value = 14510.883
$ws.Range("I126").Value = 44114.668
$ws.Range("J126").Value = 8167.2144
$ws.Range("K126").Value = 132344.004
$ws.Range("L126").Value = 24501.6432
$ws.Range("M126").Value = -129874.004
$ws.Range("N126").Value = -29441.6432
$ws.Range("H132").Value = 3730.8572
$ws.Range("I132").Value = 3676.2104
$ws.Range("K132").Value = 11028.6312
$ws.Range("M132").Value = -8498.6312
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H132").Value = 1494425.5
$ws.Range("I132").Value = 1866907
$ws.Range("K132").Value = 5600721
$ws.Range("M132").Value = -5598191
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 11984.5
$ws.Range("I81").Value = 15265
$ws.Range("J81").Value = 4330
$ws.Range("K81").Value = 30530
$ws.Range("L81").Value = 8660
$ws.Range("M81").Value = -29469
$ws.Range("N81").Value = -10782
$ws.Range("H84").Value = 11984.5
$ws.Range("I84").Value = 15265
$ws.Range("J84").Value = 4330
$ws.Range("K84").Value = 152650
$ws.Range("L84").Value = 43300
$ws.Range("M84").Value = -147346
$ws.Range("N84").Value = -53908
$ws.Range("H101").Value = 12333.444
$ws.Range("J101").Value = 12333.444
$ws.Range("L101").Value = 12333.444
$ws.Range("N101").Value = -18823.444
$ws.Range("H122").Value = 5729.3423
$ws.Range("I122").Value = 2944.4375
$ws.Range("J122").Value = 7754.727
$ws.Range("K122").Value = 8833.3125
$ws.Range("L122").Value = 23264.181
$ws.Range("M122").Value = -6383.3125
$ws.Range("N122").Value = -28164.181
$ws.Range("H132").Value = 20538.408
$ws.Range("I132").Value = 31896.154
$ws.Range("K132").Value = 95688.462
$ws.Range("M132").Value = -93158.462
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
